$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the stray "/100" suffix left over from the old formula-text convention;
# acs_data needs to be rerun before re-aggregating (per commit message), but the
# shared-string labels should just carry the bare variable name now.
$ws.Range("C41").Value = "DP05_0027E"
$ws.Range("C42").Value = "DP05_0031E"
$ws.Range("C39").Value = "DP05_0003E"

# Resize columns A and B slightly (values chosen so the persisted OOXML
# <col> width lands on 17.5 / 15.5 after the engine's width->XML padding).
$ws.Columns.Item(1).ColumnWidth = 16.666666666666668
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# Move the selection down to where the user was working.
$ws.Range("E19").Select()
